$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.52
$ws.Range("G2").Value = 1.93
$ws.Range("K2").Value = 7.2
$ws.Range("N2").Value = 2.26
$ws.Range("P2").Value = 1.51
$ws.Range("Q2").Value = 1.7
$ws.Range("T2").Value = 1.05
$ws.Range("U2").Value = 2.06
$ws.Range("F3").Value = 2.38
$ws.Range("G3").Value = 2.4
$ws.Range("L3").Value = 1.38
$ws.Range("N3").Value = 3.45
$ws.Range("Q3").Value = 2.22
$ws.Range("R3").Value = 1.87
$ws.Range("U3").Value = 1.71
$ws.Range("Z3").Value = 9.199999999999999
$ws.Range("AA4").Value = 44
$ws.Range("AB4").Value = 130
$ws.Range("AH4").Value = 7.8
$ws.Range("J4").Value = 18
$ws.Range("K4").Value = 18.5
$ws.Range("L4").Value = 1.14
$ws.Range("Q4").Value = 1.3
$ws.Range("R4").Value = 2.58
$ws.Range("S4").Value = 1.57
$ws.Range("U4").Value = 12
$ws.Range("X4").Value = 470
$ws.Range("P5").Value = 2.2
$ws.Range("W5").Value = 25
$ws.Range("G6").Value = 5.3
$ws.Range("I6").Value = 2.02
$ws.Range("J6").Value = 3.3
$ws.Range("N6").Value = 2.76
$ws.Range("Q6").Value = 1.66
$ws.Range("T6").Value = 1.98
$ws.Range("AJ7").Value = 44
$ws.Range("AK7").Value = 130
$ws.Range("AM7").Value = 70
$ws.Range("F7").Value = 2.06
$ws.Range("G7").Value = 2.08
$ws.Range("M7").Value = 1.08
$ws.Range("P7").Value = 1.81
$ws.Range("U7").Value = 1.92
$ws.Range("S8").Value = 2.22
$ws.Range("AM9").Value = 46
$ws.Range("F9").Value = 2.98
$ws.Range("G9").Value = 3
$ws.Range("H9").Value = 2.8
$ws.Range("I9").Value = 2.82
$ws.Range("N9").Value = 2.74
$ws.Range("O9").Value = 1.54
$ws.Range("Q9").Value = 2.6
$ws.Range("R9").Value = 2.14
$ws.Range("S9").Value = 1.8
$ws.Range("T9").Value = 1.55
$ws.Range("V9").Value = 8.6
$ws.Range("Z9").Value = 8.800000000000001
$ws.Range("J10").Value = 3.3
$ws.Range("K10").Value = 3.35
$ws.Range("Q10").Value = 2.3
$ws.Range("S10").Value = 1.93
$ws.Range("AC11").Value = 200
$ws.Range("AH11").Value = 10
$ws.Range("F11").Value = 1.31
$ws.Range("N11").Value = 4.8
$ws.Range("P11").Value = 2.4
$ws.Range("X11").Value = 120
$ws.Range("Y11").Value = 530
$ws.Range("I12").Value = 13.5
$ws.Range("L12").Value = 1.27
$ws.Range("P12").Value = 2.8
$ws.Range("S12").Value = 1.9
$ws.Range("Y12").Value = 490
$ws.Range("S13").Value = 2.1
$ws.Range("AA14").Value = 10.5
$ws.Range("Q14").Value = 1.45
$ws.Range("W14").Value = 18
$ws.Range("AB15").Value = 14.5
$ws.Range("AC15").Value = 44
$ws.Range("AG15").Value = 60
$ws.Range("AM15").Value = 46
$ws.Range("F15").Value = 2.4
$ws.Range("I15").Value = 3.45
$ws.Range("L15").Value = 1.38
$ws.Range("O15").Value = 1.37
$ws.Range("T15").Value = 1.4
$ws.Range("V15").Value = 11.5
$ws.Range("Y15").Value = 65
$ws.Range("AD16").Value = 11.5
$ws.Range("AH16").Value = 12.5
$ws.Range("AK16").Value = 75
$ws.Range("AL16").Value = 3.7
$ws.Range("R16").Value = 1.66
$ws.Range("Z16").Value = 15.5
$ws.Range("AA17").Value = 11.5
$ws.Range("AB17").Value = 30
$ws.Range("G17").Value = 1.44
$ws.Range("J17").Value = 5.3
$ws.Range("N17").Value = 5.1
$ws.Range("P18").Value = 1.87
$ws.Range("Q18").Value = 2.12
$ws.Range("X18").Value = 12
$ws.Range("F19").Value = 1.28
$ws.Range("N19").Value = 2.08
$ws.Range("P19").Value = 2.08
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 2.88
$ws.Range("H20").Value = 2.66
$ws.Range("I20").Value = 4.3
$ws.Range("J20").Value = 2.88
$ws.Range("N20").Value = 2.52
$ws.Range("P20").Value = 1.76
$ws.Range("Q20").Value = 1.71
$ws.Range("T20").Value = 1.3
$ws.Range("U20").Value = 1.53
